# Add data for 2024-06-10: update the 2023 (J) and 2024 (K) year-to-date
# crime-count columns across the citywide summary, the by-neighborhood
# rollup, and every individual neighborhood sheet that was affected.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3361
$ws.Range("K3").Value = 3331
$ws.Range("J4").Value = 1819
$ws.Range("K4").Value = 697
$ws.Range("K5").Value = 220
$ws.Range("K6").Value = 3912
$ws.Range("J7").Value = 29290
$ws.Range("K7").Value = 11521

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 34
$ws.Range("K3").Value = 32
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 221
$ws.Range("K4").Value = 42
$ws.Range("K7").Value = 759

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K2").Value = 90
$ws.Range("J4").Value = 25
$ws.Range("J7").Value = 593
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 127
$ws.Range("K3").Value = 171
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 460

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 66
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 109
$ws.Range("K3").Value = 139
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 402

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 84
$ws.Range("K3").Value = 71
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 57
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 201

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 11
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 40
$ws.Range("K7").Value = 325
$ws.Range("K8").Value = 759
$ws.Range("K12").Value = 19
$ws.Range("K14").Value = 59
$ws.Range("K19").Value = 352
$ws.Range("K20").Value = 264
$ws.Range("K23").Value = 111
$ws.Range("K29").Value = 604
$ws.Range("K30").Value = 38
$ws.Range("K32").Value = 18
$ws.Range("K33").Value = 460
$ws.Range("K34").Value = 55
$ws.Range("K36").Value = 135
$ws.Range("K37").Value = 402
$ws.Range("K41").Value = 99
$ws.Range("K42").Value = 412
$ws.Range("K44").Value = 106
$ws.Range("K45").Value = 13
$ws.Range("K48").Value = 143
$ws.Range("K49").Value = 68
$ws.Range("K50").Value = 67
$ws.Range("K51").Value = 136
$ws.Range("K52").Value = 314
$ws.Range("K53").Value = 152
$ws.Range("K54").Value = 225
$ws.Range("K55").Value = 122
$ws.Range("K63").Value = 38
$ws.Range("K65").Value = 275
$ws.Range("K66").Value = 43
$ws.Range("K67").Value = 449
$ws.Range("K73").Value = 102
$ws.Range("K76").Value = 176
$ws.Range("K77").Value = 82
$ws.Range("J79").Value = 801
$ws.Range("K79").Value = 298
$ws.Range("J83").Value = 593
$ws.Range("K83").Value = 248
$ws.Range("K84").Value = 81
$ws.Range("K85").Value = 542
$ws.Range("K89").Value = 153
$ws.Range("K90").Value = 103
$ws.Range("K91").Value = 120
$ws.Range("K94").Value = 142
$ws.Range("K95").Value = 187
$ws.Range("K96").Value = 141
$ws.Range("K99").Value = 201
$ws.Range("K100").Value = 21
$ws.Range("J101").Value = 29290
$ws.Range("K101").Value = 11521

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 139
$ws.Range("K3").Value = 146
$ws.Range("K4").Value = 25
$ws.Range("K6").Value = 131
$ws.Range("K7").Value = 449

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 67
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 205
$ws.Range("K5").Value = 14
$ws.Range("K6").Value = 184
$ws.Range("K7").Value = 604

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 119
$ws.Range("K3").Value = 93
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 352

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K6").Value = 50
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 32
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 176

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 59

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 109
$ws.Range("K3").Value = 131
$ws.Range("K6").Value = 153
$ws.Range("K7").Value = 412

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 31
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J4").Value = 45
$ws.Range("K6").Value = 68
$ws.Range("J7").Value = 801
$ws.Range("K7").Value = 298

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 92
$ws.Range("K3").Value = 74
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K2").Value = 5
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 21

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 117
$ws.Range("K7").Value = 325

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 18
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K4").Value = 15
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("K2").Value = 12
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 31
$ws.Range("K4").Value = 6
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K3").Value = 33
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 103

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 35
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 185
$ws.Range("K7").Value = 542

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K2").Value = 40
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 82

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 13

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 83
$ws.Range("K5").Value = 10
$ws.Range("K6").Value = 125
$ws.Range("K7").Value = 314

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 19
